$d = $word.ActiveDocument

# The document currently starts with:
#   <w:p pStyle=Heading1>Day After Day - June 1937</w:p>
#   <w:p><w:r><w:rPr><w:b/></w:rPr>By Dorothy Day</w:p>
# (wrapped in a bookmark around the title paragraph).
#
# We turn this into a pandoc-style title block:
#   <w:p pStyle=Title>Day After Day - June 1937</w:p>   (split word-by-word into runs)
#   <w:p pStyle=Authors>Dorothy Day</w:p>                (split word-by-word into runs)

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)
$range = $d.Range($p1.Range.Start, $p2.Range.End)

$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Title"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">After</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">-</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">June</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">1937</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Authors"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$range.InsertXML($newXml)
